$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.406.88"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "3.180.43"
$ws.Range("E3").Value = "  +4.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.177.68"
$ws.Range("E8").Value = "  +4.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("E10").Value = "  +6.64%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.506"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.09%  "
$ws.Range("E13").Value = "  +19.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.25%  "
$ws.Range("D15").Value = "3.702.51"
$ws.Range("E15").Value = "  +4.28%  "
$ws.Range("D16").Value = "65.339.55"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.58%  "
$ws.Range("D18").Value = "3.179.01"
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "515.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.93%  "
$ws.Range("E23").Value = "  +6.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.45%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.02%  "
$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +16.43%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.997"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.96%  "
$ws.Range("E34").Value = "  +11.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0907"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "480.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.62%  "
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.81%  "
$ws.Range("D42").Value = "3.088.61"
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("E43").Value = "  +3.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.287"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.52%  "
$ws.Range("D47").Value = "0.0₃0610"
$ws.Range("E47").Value = "  +19.19%  "
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("E50").Value = "  +8.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.19%  "
